# Auto-applied numeric updates to the Leve-profit sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) reflecting refreshed Market Board price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 11634.308
$ws.Range("I125").Value = 18910.857
$ws.Range("J125").Value = 3145
$ws.Range("K125").Value = 170197.713
$ws.Range("L125").Value = 28305
$ws.Range("M125").Value = -167737.713
$ws.Range("N125").Value = -33225

$ws.Range("H138").Value = 2609.795
$ws.Range("I138").Value = 1204.409
$ws.Range("J138").Value = 4428.5293
$ws.Range("K138").Value = 3613.227
$ws.Range("L138").Value = 13285.5879
$ws.Range("M138").Value = 1526.773
$ws.Range("N138").Value = -23565.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1612.3636
$ws.Range("I2").Value = 868.3333
$ws.Range("J2").Value = 2505.2
$ws.Range("K2").Value = 868.3333
$ws.Range("L2").Value = 2505.2
$ws.Range("M2").Value = -755.3333
$ws.Range("N2").Value = -2731.2

$ws.Range("H45").Value = 1193.375
$ws.Range("I45").Value = 1149.5714
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1149.5714
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -772.5714
$ws.Range("N45").Value = -2254

$ws.Range("H112").Value = 11752.223
$ws.Range("J112").Value = 11752.223
$ws.Range("L112").Value = 11752.223
$ws.Range("N112").Value = -14706.223

$ws.Range("H116").Value = 1612.3636
$ws.Range("I116").Value = 868.3333
$ws.Range("J116").Value = 2505.2
$ws.Range("K116").Value = 868.3333
$ws.Range("L116").Value = 2505.2
$ws.Range("M116").Value = 1425.6667
$ws.Range("N116").Value = -7093.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1612.3636
$ws.Range("I3").Value = 868.3333
$ws.Range("J3").Value = 2505.2
$ws.Range("K3").Value = 868.3333
$ws.Range("L3").Value = 2505.2
$ws.Range("M3").Value = -754.3333
$ws.Range("N3").Value = -2733.2

$ws.Range("H107").Value = 1106.75
$ws.Range("I107").Value = 864.55554
$ws.Range("J107").Value = 1833.3334
$ws.Range("K107").Value = 864.55554
$ws.Range("L107").Value = 1833.3334
$ws.Range("M107").Value = 1055.44446
$ws.Range("N107").Value = -5673.3334

$ws.Range("H110").Value = 22222
$ws.Range("J110").Value = 22222
$ws.Range("L110").Value = 22222
$ws.Range("N110").Value = -30402

$ws.Range("H122").Value = 48780
$ws.Range("J122").Value = 48780
$ws.Range("L122").Value = 48780
$ws.Range("N122").Value = -58580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3686.6667
$ws.Range("I16").Value = 1998
$ws.Range("J16").Value = 5797.5
$ws.Range("K16").Value = 1998
$ws.Range("L16").Value = 5797.5
$ws.Range("M16").Value = -1711
$ws.Range("N16").Value = -6371.5

$ws.Range("H52").Value = 24840
$ws.Range("J52").Value = 24840
$ws.Range("L52").Value = 24840
$ws.Range("N52").Value = -25428

$ws.Range("H105").Value = 882.8571
$ws.Range("I105").Value = 921
$ws.Range("K105").Value = 921
$ws.Range("M105").Value = 826

$ws.Range("H107").Value = 84380.336
$ws.Range("I107").Value = 91960.37
$ws.Range("K107").Value = 91960.37
$ws.Range("M107").Value = -90040.37

$ws.Range("H113").Value = 3686.6667
$ws.Range("I113").Value = 1998
$ws.Range("J113").Value = 5797.5
$ws.Range("K113").Value = 1998
$ws.Range("L113").Value = 5797.5
$ws.Range("M113").Value = 172
$ws.Range("N113").Value = -10137.5

$ws.Range("H135").Value = 40541.75
$ws.Range("J135").Value = 40541.75
$ws.Range("L135").Value = 40541.75
$ws.Range("N135").Value = -50681.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 431.25
$ws.Range("I107").Value = 209.76315
$ws.Range("J107").Value = 1834
$ws.Range("K107").Value = 629.28945
$ws.Range("L107").Value = 5502
$ws.Range("M107").Value = 1290.71055
$ws.Range("N107").Value = -9342

$ws.Range("H118").Value = 1900.6428
$ws.Range("I118").Value = 518.1667
$ws.Range("J118").Value = 2937.5
$ws.Range("K118").Value = 1554.5001
$ws.Range("L118").Value = 8812.5
$ws.Range("M118").Value = -311.5001
$ws.Range("N118").Value = -11298.5

$ws.Range("H119").Value = 1434.8334
$ws.Range("I119").Value = 527.25
$ws.Range("K119").Value = 1581.75
$ws.Range("M119").Value = 3256.25

$ws.Range("H120").Value = 6333.3335
$ws.Range("I120").Value = 2000
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 6000
$ws.Range("L120").Value = 45000
$ws.Range("M120").Value = -1162
$ws.Range("N120").Value = -54676

$ws.Range("H121").Value = 20653.3
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 20653.3
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 61959.89999999999
$ws.Range("N121").Value = -64579.89999999999
$ws.Range("M121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3214.52
$ws.Range("I113").Value = 2973.4583
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 2973.4583
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -803.4582999999998
$ws.Range("N113").Value = -13340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1611.5
$ws.Range("I61").Value = 1492
$ws.Range("J61").Value = 1671.25
$ws.Range("K61").Value = 1492
$ws.Range("L61").Value = 1671.25
$ws.Range("M61").Value = -1290
$ws.Range("N61").Value = -2075.25

$ws.Range("H113").Value = 1611.5
$ws.Range("I113").Value = 1492
$ws.Range("J113").Value = 1671.25
$ws.Range("K113").Value = 1492
$ws.Range("L113").Value = 1671.25
$ws.Range("M113").Value = 678
$ws.Range("N113").Value = -6011.25

$ws.Range("H136").Value = 4351.838
$ws.Range("I136").Value = 1456.9565
$ws.Range("J136").Value = 9107.714
$ws.Range("K136").Value = 4370.8695
$ws.Range("L136").Value = 27323.142
$ws.Range("M136").Value = -1820.8695
$ws.Range("N136").Value = -32423.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 18367.715
$ws.Range("J112").Value = 18367.715
$ws.Range("L112").Value = 18367.715
$ws.Range("N112").Value = -21321.715

$ws.Range("H113").Value = 101074.5
$ws.Range("I113").Value = 201359.4
$ws.Range("J113").Value = 789.6
$ws.Range("K113").Value = 604078.2
$ws.Range("L113").Value = 2368.8
$ws.Range("M113").Value = -601908.2
$ws.Range("N113").Value = -6708.8

$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800

$ws.Range("H126").Value = 1029.3684
$ws.Range("I126").Value = 899.5
$ws.Range("J126").Value = 1393
$ws.Range("K126").Value = 2698.5
$ws.Range("L126").Value = 4179
$ws.Range("M126").Value = -228.5
$ws.Range("N126").Value = -9119

$ws.Range("H136").Value = 3993.8635
$ws.Range("I136").Value = 5718.864
$ws.Range("K136").Value = 17156.592
$ws.Range("M136").Value = -14606.592

